# Record new S5/S6 (columns F/G) scores for the "Mago" sheet, then make
# "Mago" the active sheet/tab (matches the sheetView tabSelected move from
# Ninja -> Mago and the workbook-level activeTab change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mago")

# David Tsegave
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0

# Juan Pablo Alfaya
$ws.Range("G3").Value = 0

# Alexandre Martínez
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0

# Ruth de la Fuente
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Alexia Torres
$ws.Range("G6").Value = 0

# Javier Janeiro
$ws.Range("G7").Value = 0

# Make "Mago" the active/selected sheet (was "Ninja").
$ws.Activate()
